$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new training-day columns (CU, CV) after the existing CT column ---
# Copy the format of the CT column's date header (CT1) and data cells into
# CU/CV so the new columns inherit the same styles as the rest of the sheet.
# (Rows 12, 21 and 23 don't have a CT entry -- those players' data stops
# earlier -- so we copy in separate blocks that skip them.)

$ws.Range("CT1:CT11").Copy()
$ws.Range("CU1:CV11").PasteSpecial(-4122)

$ws.Range("CT13:CT20").Copy()
$ws.Range("CU13:CV20").PasteSpecial(-4122)

$ws.Range("CT22").Copy()
$ws.Range("CU22:CV22").PasteSpecial(-4122)

$ws.Range("CT24:CT30").Copy()
$ws.Range("CU24:CV30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New session dates (row 1 header)
$ws.Cells.Item(1, 99).Value = 46027
$ws.Cells.Item(1, 100).Value = 46028

# Attendance marks for the two new sessions, per player row
$ws.Cells.Item(2, 99).Value = "P"
$ws.Cells.Item(2, 100).Value = "P"
$ws.Cells.Item(3, 99).Value = "P"
$ws.Cells.Item(3, 100).Value = "P"
$ws.Cells.Item(4, 99).Value = "P"
$ws.Cells.Item(4, 100).Value = "P"
$ws.Cells.Item(5, 99).Value = "P"
$ws.Cells.Item(5, 100).Value = "P"
$ws.Cells.Item(6, 99).Value = "A"
$ws.Cells.Item(6, 100).Value = "P"
$ws.Cells.Item(7, 99).Value = "P"
$ws.Cells.Item(7, 100).Value = "P"
$ws.Cells.Item(8, 99).Value = "P"
$ws.Cells.Item(8, 100).Value = "P"
$ws.Cells.Item(9, 99).Value = "M"
$ws.Cells.Item(9, 100).Value = "M"
$ws.Cells.Item(10, 99).Value = "P"
$ws.Cells.Item(10, 100).Value = "P"
$ws.Cells.Item(11, 99).Value = "P"
$ws.Cells.Item(11, 100).Value = "P"
$ws.Cells.Item(13, 99).Value = "B"
$ws.Cells.Item(13, 100).Value = "B"
$ws.Cells.Item(14, 99).Value = "P"
$ws.Cells.Item(14, 100).Value = "P"
$ws.Cells.Item(15, 99).Value = "P"
$ws.Cells.Item(15, 100).Value = "P"
$ws.Cells.Item(16, 99).Value = "RH"
$ws.Cells.Item(16, 100).Value = "P"
$ws.Cells.Item(17, 99).Value = "P"
$ws.Cells.Item(17, 100).Value = "P"
$ws.Cells.Item(18, 99).Value = "P"
$ws.Cells.Item(18, 100).Value = "P"
$ws.Cells.Item(19, 99).Value = "P"
$ws.Cells.Item(19, 100).Value = "P"
$ws.Cells.Item(20, 99).Value = "P"
$ws.Cells.Item(20, 100).Value = "P"
$ws.Cells.Item(22, 99).Value = "P"
$ws.Cells.Item(22, 100).Value = "P"
$ws.Cells.Item(24, 99).Value = "P"
$ws.Cells.Item(24, 100).Value = "P"
$ws.Cells.Item(25, 99).Value = "A"
$ws.Cells.Item(25, 100).Value = "A"
$ws.Cells.Item(26, 99).Value = "M"
$ws.Cells.Item(26, 100).Value = "M"
$ws.Cells.Item(27, 99).Value = "P"
$ws.Cells.Item(27, 100).Value = "P"
$ws.Cells.Item(28, 99).Value = "P"
$ws.Cells.Item(28, 100).Value = "P"
$ws.Cells.Item(29, 99).Value = "P"
$ws.Cells.Item(29, 100).Value = "P"
$ws.Cells.Item(30, 99).Value = "P"
$ws.Cells.Item(30, 100).Value = "P"

# --- Update the active selection to reflect the new right-most data entry point ---
# (the sheet is already frozen at column A / xSplit=1; we keep that and just
# move the cursor, matching the author scrolling over to work on the new columns)
$ws.Range("CX27").Select()
